$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exportar")

# Update Data de Inicio / Data Final columns (G2:G19, H2:H19)
$ws.Range("G2:G19").Value = "14/07/2025 00:00"
$ws.Range("H2:H19").Value = "14/07/2025 23:59"

# Update numeric metric columns per diff
$ws.Range("I2").Value = 11.75
$ws.Range("J2").Value = 8640.9
$ws.Range("K2").Value = 43.144257222222215
$ws.Range("L2").Value = 128.25583749999998
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 6.485298931423896
$ws.Range("O2").Value = 20.791240116174354
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 171.40009583333335
$ws.Range("R2").Value = 565.5364140101952
$ws.Range("S2").Value = 0.588
$ws.Range("T2").Value = 13.403557774303266
$ws.Range("U2").Value = 44.22517946564514
$ws.Range("V2").Value = 0.2735561070501127
$ws.Range("J4").Value = 8254.45
$ws.Range("I7").Value = 12.949999999998909
$ws.Range("J7").Value = 8891
$ws.Range("K7").Value = 37.514920277777776
$ws.Range("L7").Value = 128.04666194444442
$ws.Range("M7").Value = 1.7534527777777777
$ws.Range("N7").Value = 6.396066852968526
$ws.Range("O7").Value = 15.766277080704343
$ws.Range("P7").Value = 20.347662779397474
$ws.Range("Q7").Value = 167.31503611111114
$ws.Range("R7").Value = 552.0577166087103
$ws.Range("S7").Value = 0.6679999999999999
$ws.Range("T7").Value = 11.916687393131028
$ws.Range("U7").Value = 39.31923504724892
$ws.Range("V7").Value = 0.2869238584198225
$ws.Range("I8").Value = 19.849999999999454
$ws.Range("J8").Value = 6938.6
$ws.Range("K8").Value = 49.53759805555555
$ws.Range("L8").Value = 179.05847416666666
$ws.Range("M8").Value = 15.846829722222223
$ws.Range("N8").Value = 4.696134463129835
$ws.Range("O8").Value = 18.0240630281001
$ws.Range("P8").Value = 31.132404060248856
$ws.Range("Q8").Value = 244.44290305555555
$ws.Range("R8").Value = 806.5419225827291
$ws.Range("S8").Value = 0.348
$ws.Range("T8").Value = 11.673197163243074
$ws.Range("U8").Value = 38.51583647977524
$ws.Range("V8").Value = 0.26152155813881867
$ws.Range("I9").Value = 20.25
$ws.Range("J9").Value = 7751.15
$ws.Range("K9").Value = 44.40174249999997
$ws.Range("L9").Value = 191.86119472222222
$ws.Range("M9").Value = 16.12675638888889
$ws.Range("N9").Value = 4.4769097016299195
$ws.Range("O9").Value = 18.142762121514572
$ws.Range("P9").Value = 32.5412461038233
$ws.Range("Q9").Value = 252.38969444444444
$ws.Range("R9").Value = 832.7624441239154
$ws.Range("S9").Value = 0.588
$ws.Range("T9").Value = 12.065113772216096
$ws.Range("U9").Value = 39.8089694418779
$ws.Range("V9").Value = 0.2624503201294111
$ws.Range("I10").Value = 12.199999999998909
$ws.Range("J10").Value = 9340.15
$ws.Range("K10").Value = 22.452311388888887
$ws.Range("L10").Value = 125.8994411111111
$ws.Range("M10").Value = 6.5919847222222225
$ws.Range("N10").Value = 4.366760849644152
$ws.Range("O10").Value = 15.558125661493689
$ws.Range("P10").Value = 29.57345411352838
$ws.Range("Q10").Value = 154.94373611111112
$ws.Range("R10").Value = 511.2384824966854
$ws.Range("S10").Value = 0.72
$ws.Range("T10").Value = 11.556279717622635
$ws.Range("U10").Value = 38.13006614160862
$ws.Range("V10").Value = 0.25883073164200204
$ws.Range("I11").Value = 10.799999999999272
$ws.Range("J11").Value = 8720.1
$ws.Range("K11").Value = 41.85617972222224
$ws.Range("L11").Value = 85.46789499999998
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 6.290916353399304
$ws.Range("O11").Value = 18.877239733002856
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 127.32407527777781
$ws.Range("R11").Value = 420.1071218756871
$ws.Range("S11").Value = 0.728
$ws.Range("T11").Value = 11.416189176630477
$ws.Range("U11").Value = 37.667835932201434
$ws.Range("V11").Value = 0.29207950126470683
$ws.Range("I12").Value = 20.34999999999991
$ws.Range("J12").Value = 1297.7
$ws.Range("K12").Value = 97.21722555555556
$ws.Range("L12").Value = 243.83604527777783
$ws.Range("M12").Value = 0.15145416666666667
$ws.Range("N12").Value = 10.36700631747301
$ws.Range("O12").Value = 20.515931265847936
$ws.Range("P12").Value = 17.08
$ws.Range("Q12").Value = 341.2047233333334
$ws.Range("R12").Value = 1125.8085635197597
$ws.Range("S12").Value = 0.54
$ws.Range("T12").Value = 16.040236089620798
$ws.Range("U12").Value = 52.92492722303892
$ws.Range("V12").Value = 0.3470115614232021
$ws.Range("I13").Value = 12.549999999999955
$ws.Range("J13").Value = 1285
$ws.Range("K13").Value = 34.99618527777777
$ws.Range("L13").Value = 128.09865416666668
$ws.Range("M13").Value = 1.5333794444444444
$ws.Range("N13").Value = 6.720012158897594
$ws.Range("O13").Value = 16.765228576102366
$ws.Range("P13").Value = 30.068526488257785
$ws.Range("Q13").Value = 164.62822027777779
$ws.Range("R13").Value = 543.192539584733
$ws.Range("S13").Value = 0.536
$ws.Range("T13").Value = 12.762036086625125
$ws.Range("U13").Value = 42.10847192825818
$ws.Range("V13").Value = 0.31177736835312886
$ws.Range("I19").Value = 14.600000000000136
$ws.Range("J19").Value = 1383.45
$ws.Range("K19").Value = 81.35426861111114
$ws.Range("L19").Value = 135.13235194444445
$ws.Range("M19").Value = 1.0984266666666667
$ws.Range("N19").Value = 10.375846784363175
$ws.Range("O19").Value = 18.8633921143646
$ws.Range("P19").Value = 24.31918411330049
$ws.Range("Q19").Value = 217.58504666666656
$ws.Range("R19").Value = 717.9241437167084
$ws.Range("S19").Value = 0.532
$ws.Range("T19").Value = 14.457798059958218
$ws.Range("U19").Value = 47.70365634604392
$ws.Range("V19").Value = 0.39373101772039976
